$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: drop the row-number / username cells (A3, B3) -- the remaining
# row-3 data (password/email/rolecode/roletype) stays as-is.
$ws.Range("A3:B3").ClearContents()

# Row 4: new data row. (Shared-string append order matters: "ADMIN, ADMIN"
# must land before "test3" to match the target string table ordering.)
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = "passjdfdfdfdff"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:test2@gmail.com")
$ws.Range("D4").Value = "test2@gmail.com"
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("F4").Value = "ADMIN, ADMIN"
$ws.Range("B4").Value = "test3"

# Selection moved.
$ws.Range("C9").Select()
